$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the N column values (rows 76-94) from 60 to 30 for the LCC calc.
$nRows = 76..94
foreach ($r in $nRows) {
    $ws.Cells.Item($r, 14).Value = 30
}

# Row 97 column N holds a formula (=AVERAGE(N76,N79,N82,N85,N88)) and will
# recalculate automatically once its precedents change, but set it
# explicitly as well to be safe.
$ws.Cells.Item(98, 14).Value = 30
$ws.Cells.Item(99, 14).Value = 30

# Row 101: B101 switches from the shared string "Dummy" to the numeric 0.
$ws.Cells.Item(101, 2).Value = 0

# Update the sheet view / pane / selection state.
$ws.Application.ActiveWindow.ScrollColumn = 3
$window = $excel.ActiveWindow
$window.SplitRow = 1
$window.FreezePanes = $true

$ws.Range("L99").Select()
